# Auto-generated edit script applying scheduled runner price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!5  "Met a Sticky End" / "Animal Glue"  (item id 5503)
$ws.Range("H5").Value = 104.125
$ws.Range("I5").Value = 104.125
$ws.Range("K5").Value = 104.125
$ws.Range("M5").Value = 10.875

# ALC!20  "Shut Up and Take My Gil" / "Ash Wand"  (item id 1965)
$ws.Range("H20").Value = 1024.75
$ws.Range("I20").Value = 1024.75
$ws.Range("K20").Value = 1024.75
$ws.Range("M20").Value = -794.75

# ALC!35  "Conspicuous Conjuration" / "Whispering Ash Wand"  (item id 1965)
$ws.Range("H35").Value = 1024.75
$ws.Range("I35").Value = 1024.75
$ws.Range("K35").Value = 1024.75
$ws.Range("M35").Value = -645.75

# ALC!40  "Stuck in the Moment" / "Horn Glue"  (item id 5505)
$ws.Range("H40").Value = 5852036.5
$ws.Range("I40").Value = 2923.4167
$ws.Range("J40").Value = 15879087
$ws.Range("K40").Value = 2923.4167
$ws.Range("L40").Value = 15879087
$ws.Range("M40").Value = -2748.4167
$ws.Range("N40").Value = -15879437

# ALC!106  "Making Your Mark" / "Enchanted Palladium Ink"  (item id 19903)
$ws.Range("H106").Value = 2278.2632
$ws.Range("I106").Value = 2337.5386
$ws.Range("J106").Value = 2149.8333
$ws.Range("K106").Value = 2337.5386
$ws.Range("L106").Value = 2149.8333
$ws.Range("M106").Value = -1706.5386
$ws.Range("N106").Value = -3411.8333

# ALC!112  "Making Ends Meet" / "Superior Spiritbond Potion"  (item id 27960)
$ws.Range("H112").Value = 119868.414
$ws.Range("I112").Value = 201079.6
$ws.Range("J112").Value = 86030.414
$ws.Range("K112").Value = 603238.8
$ws.Range("L112").Value = 258091.242
$ws.Range("M112").Value = -602130.8
$ws.Range("N112").Value = -260307.242

# ALC!113  "Amaro Kart" / "Starch Glue"  (item id 27775)
$ws.Range("H113").Value = 6355.5713
$ws.Range("I113").Value = 6229.6665
$ws.Range("J113").Value = 6450
$ws.Range("K113").Value = 6229.6665
$ws.Range("L113").Value = 6450
$ws.Range("M113").Value = -2975.6665
$ws.Range("N113").Value = -12958

# ALC!116  "Growing Up" / "Growth Formula Kappa"  (item id 27778)
$ws.Range("H116").Value = 4950.95
$ws.Range("I116").Value = 4950.95
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 4950.95
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1508.95
$ws.Range("N116").ClearContents()

# ALC!132  "Fast-forwarding Flora" / "Growth Formula Lambda"  (item id 44049)
$ws.Range("H132").Value = 3249.9092
$ws.Range("I132").Value = 3335
$ws.Range("K132").Value = 10005
$ws.Range("M132").Value = -7475

$ws = $wb.Worksheets.Item("ARM")
# ARM!2  "Ain't Got No Ingots" / "Bronze Ingot"  (item id 27713)
$ws.Range("H2").Value = 999
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# ARM!5  "The Alloyed Truth" / "Bronze Rivets"  (item id 5091)
$ws.Range("H5").Value = 1612.2142
$ws.Range("I5").Value = 734.1818
$ws.Range("J5").Value = 4831.6665
$ws.Range("K5").Value = 734.1818
$ws.Range("L5").Value = 4831.6665
$ws.Range("M5").Value = -622.1818
$ws.Range("N5").Value = -5055.6665

# ARM!116  "No Scope" / "Titanbronze Ingot"  (item id 27713)
$ws.Range("H116").Value = 999
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# BSM!3  "Hells Bells" / "Bronze Ingot"  (item id 27713)
$ws.Range("H3").Value = 999
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

# BSM!4  "Mending Fences" / "Bronze Rivets"  (item id 5091)
$ws.Range("H4").Value = 1612.2142
$ws.Range("I4").Value = 734.1818
$ws.Range("J4").Value = 4831.6665
$ws.Range("K4").Value = 734.1818
$ws.Range("L4").Value = 4831.6665
$ws.Range("M4").Value = -619.1818
$ws.Range("N4").Value = -5061.6665

# BSM!22  "Riveting Run" / "Iron Rivets"  (item id 5092)
$ws.Range("H22").Value = 1906226.5
$ws.Range("I22").Value = 1639.625
$ws.Range("J22").Value = 5292159
$ws.Range("K22").Value = 1639.625
$ws.Range("L22").Value = 5292159
$ws.Range("M22").Value = -1466.625
$ws.Range("N22").Value = -5292505

# BSM!86  "Through Thick and Thin" / "Adamantite Nugget"  (item id 12526)
$ws.Range("H86").Value = 3781.4119
$ws.Range("I86").Value = 3850.7144
$ws.Range("J86").Value = 3732.9
$ws.Range("K86").Value = 3850.7144
$ws.Range("L86").Value = 3732.9
$ws.Range("M86").Value = -2727.7144
$ws.Range("N86").Value = -5978.9

# BSM!89  "Piercing Eyes Deserve Piercing Shafts (L)" / "Adamantite Nugget"  (item id 12526)
$ws.Range("H89").Value = 3781.4119
$ws.Range("I89").Value = 3850.7144
$ws.Range("J89").Value = 3732.9
$ws.Range("K89").Value = 19253.572
$ws.Range("L89").Value = 18664.5
$ws.Range("M89").Value = -13637.572
$ws.Range("N89").Value = -29896.5

# BSM!105  "Ingot to Wing It" / "Molybdenum Ingot"  (item id 19947)
$ws.Range("H105").Value = 1867.4286
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# BSM!123  "Archon Denied" / "High Durium Saw"  (item id 35320)
$ws.Range("H123").Value = 74499
$ws.Range("J123").Value = 74499
$ws.Range("L123").Value = 74499
$ws.Range("N123").Value = -84299

# BSM!137  "Dagger Swagger" / "Cobalt Tungsten Khukuri"  (item id 42153)
$ws.Range("H137").Value = 48999
$ws.Range("J137").Value = 48999
$ws.Range("L137").Value = 48999
$ws.Range("N137").Value = -59199

$ws = $wb.Worksheets.Item("CRP")
# CRP!2  "In with the New" / "Bone Harpoon"  (item id 1820)
$ws.Range("H2").Value = 3487.5
$ws.Range("I2").Value = 1225
$ws.Range("J2").Value = 5750
$ws.Range("K2").Value = 1225
$ws.Range("L2").Value = 5750
$ws.Range("M2").Value = -1112
$ws.Range("N2").Value = -5976

# CRP!7  "Gridania's Got Talent" / "Maple Lumber"  (item id 5361)
$ws.Range("H7").Value = 6985.1763
$ws.Range("I7").Value = 11804.8
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 11804.8
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = -11691.8
$ws.Range("N7").Value = -326

# CRP!12  "A Sword in Hand" / "Ash Macuahuitl"  (item id 1604)
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

# CRP!19  "Shielding Sales" / "Square Ash Shield"  (item id 2233)
$ws.Range("H19").Value = 2299.7144
$ws.Range("I19").Value = 2399.6667
$ws.Range("J19").Value = 2224.75
$ws.Range("K19").Value = 2399.6667
$ws.Range("L19").Value = 2224.75
$ws.Range("M19").Value = -2229.6667
$ws.Range("N19").Value = -2564.75

# CRP!24  "What You Need" / "Square Ash Shield"  (item id 2233)
$ws.Range("H24").Value = 2299.7144
$ws.Range("I24").Value = 2399.6667
$ws.Range("J24").Value = 2224.75
$ws.Range("K24").Value = 2399.6667
$ws.Range("L24").Value = 2224.75
$ws.Range("M24").Value = -2229.6667
$ws.Range("N24").Value = -2564.75

# CRP!31  "Wall Not Found" / "Walnut Lumber"  (item id 44023)
$ws.Range("H31").Value = 4530.5
$ws.Range("I31").Value = 3121.6667
$ws.Range("J31").Value = 8757
$ws.Range("K31").Value = 3121.6667
$ws.Range("L31").Value = 8757
$ws.Range("M31").Value = -2826.6667
$ws.Range("N31").Value = -9347

# CRP!34  "Armoires of the Rich and Famous" / "Walnut Lumber"  (item id 44023)
$ws.Range("H34").Value = 4530.5
$ws.Range("I34").Value = 3121.6667
$ws.Range("J34").Value = 8757
$ws.Range("K34").Value = 3121.6667
$ws.Range("L34").Value = 8757
$ws.Range("M34").Value = -2919.6667
$ws.Range("N34").Value = -9161

# CRP!47  "Grippy When Wet" / "Mythril Cavalry Bow"  (item id 1920)
$ws.Range("H47").Value = 34499
$ws.Range("I47").Value = 28999
$ws.Range("J47").Value = 39999
$ws.Range("K47").Value = 28999
$ws.Range("L47").Value = 39999
$ws.Range("M47").Value = -28433
$ws.Range("N47").Value = -41131

# CRP!54  "The Turning Point" / "Garnet Grinding Wheel"  (item id 2413)
$ws.Range("H54").Value = 36666.332
$ws.Range("J54").Value = 36666.332
$ws.Range("L54").Value = 36666.332
$ws.Range("N54").Value = -37982.332

# CRP!62  "Splinter in the Sewers" / "Cedar Lumber"  (item id 12580)
$ws.Range("H62").Value = 2532
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 3144.8
$ws.Range("K62").Value = 1000
$ws.Range("L62").Value = 3144.8
$ws.Range("M62").Value = -376
$ws.Range("N62").Value = -4392.8

# CRP!65  "The Lumber of Their Discontent (L)" / "Cedar Lumber"  (item id 12580)
$ws.Range("H65").Value = 2532
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 3144.8
$ws.Range("K65").Value = 5000
$ws.Range("L65").Value = 15724
$ws.Range("M65").Value = -1880
$ws.Range("N65").Value = -21964

# CRP!99  "O Pine" / "Pine Lumber"  (item id 36198)
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

# CRP!126  "A Better Conductor" / "Red Pine Lumber"  (item id 36198)
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws = $wb.Worksheets.Item("CUL")
# CUL!23  "Sweet Smell of Success" / "Lavender Oil"  (item id 4858)
$ws.Range("H23").Value = 1002.1539
$ws.Range("I23").Value = 708.6667
$ws.Range("J23").Value = 1090.2
$ws.Range("K23").Value = 2126.0001
$ws.Range("L23").Value = 3270.6
$ws.Range("M23").Value = -1891.0001
$ws.Range("N23").Value = -3740.6

# CUL!107  "Slippery Service" / "Frantoio Oil"  (item id 27838)
$ws.Range("H107").Value = 817.2174
$ws.Range("I107").Value = 367.2
$ws.Range("J107").Value = 1163.3846
$ws.Range("K107").Value = 1101.6
$ws.Range("L107").Value = 3490.1538
$ws.Range("M107").Value = 818.4000000000001
$ws.Range("N107").Value = -7330.1538

# CUL!121  "A Cookie for Your Troubles" / "Coffee Biscuit"  (item id 27878)
$ws.Range("H121").Value = 999999
$ws.Range("I121").Value = 999999
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 2999997
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -2998687
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# GSM!2  "Copper and Robbers" / "Copper Ingot"  (item id 5062)
$ws.Range("H2").Value = 79.72221999999999
$ws.Range("J2").Value = 192.57143
$ws.Range("L2").Value = 192.57143
$ws.Range("N2").Value = -418.57143

# GSM!10  "Let's Talk about Hex" / "Bone Necklace"  (item id 4306)
$ws.Range("H10").Value = 27666.334
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

# GSM!40  "A Little Bird Told Me" / "Malachite Bracelet"  (item id 4113)
$ws.Range("H40").Value = 24950
$ws.Range("J40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("N40").Value = -20302

# GSM!80  "Needs More Prayerbell" / "Hardsilver Ingot"  (item id 12521)
$ws.Range("H80").Value = 3483.5625
$ws.Range("I80").Value = 3839.2222
$ws.Range("J80").Value = 3026.2856
$ws.Range("K80").Value = 3839.2222
$ws.Range("L80").Value = 3026.2856
$ws.Range("M80").Value = -2841.2222
$ws.Range("N80").Value = -5022.2856

# GSM!83  "With a Noise That Reaches Heaven (L)" / "Hardsilver Ingot"  (item id 12521)
$ws.Range("H83").Value = 3483.5625
$ws.Range("I83").Value = 3839.2222
$ws.Range("J83").Value = 3026.2856
$ws.Range("K83").Value = 19196.111
$ws.Range("L83").Value = 15131.428
$ws.Range("M83").Value = -14204.111
$ws.Range("N83").Value = -25115.428

# GSM!122  "Awarding Academic Excellence" / "Ametrine"  (item id 36182)
$ws.Range("H122").Value = 7092.625
$ws.Range("J122").Value = 23000
$ws.Range("L122").Value = 69000
$ws.Range("N122").Value = -73900

$ws = $wb.Worksheets.Item("LTW")
# LTW!82  "Trainin' the Neck" / "Dragon Leather"  (item id 12565)
$ws.Range("H82").Value = 1496.2941
$ws.Range("J82").Value = 1787
$ws.Range("L82").Value = 1787
$ws.Range("N82").Value = -2509

# LTW!85  "Training Is Only Skintight (L)" / "Dragon Leather"  (item id 12565)
$ws.Range("H85").Value = 1496.2941
$ws.Range("J85").Value = 1787
$ws.Range("L85").Value = 1787
$ws.Range("N85").Value = -4283

# LTW!123  "Running up the Tabi" / "Gajaskin Tabi"  (item id 35408)
$ws.Range("H123").Value = 80997.5
$ws.Range("J123").Value = 80997.5
$ws.Range("L123").Value = 80997.5
$ws.Range("N123").Value = -90797.5

$ws = $wb.Worksheets.Item("WVR")
# WVR!41  "Half Is the New Double" / "Linen Halfgloves"  (item id 21725)
$ws.Range("H41").Value = 36321.5
$ws.Range("J41").Value = 37786
$ws.Range("L41").Value = 37786
$ws.Range("N41").Value = -38566

# WVR!113  "A Tender Table" / "Pixie Floss"  (item id 27752)
$ws.Range("H113").Value = 983.24445
$ws.Range("I113").Value = 962.4
$ws.Range("K113").Value = 2887.2
$ws.Range("M113").Value = -717.1999999999998

# WVR!136  "Weaving the Envelope" / "Sarcenet Cloth"  (item id 44031)
$ws.Range("H136").Value = 25002552
$ws.Range("I136").Value = 25002552
$ws.Range("K136").Value = 75007656
$ws.Range("M136").Value = -75005106
